# Update cryptocurrency price (column D) and 1h volume change (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.403.20"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "'1.850.67"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'240.74"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'0.6293"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.07674"

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").Value = "'24.55"
$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").Value = "'0.07748"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "'1.856.98"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "'0.00001113"
$ws.Range("E13").Value = "  +10.46%  "

$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").Value = "'0.6807"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("D16").Value = "'83.59"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "'2.108.22"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "'6.147"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'29.463.84"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'229.13"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "'7.442"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("D24").Value = "'1.0000"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'156.89"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'0.1384"
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").Value = "'8.383"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").Value = "'17.67"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  +3.74%  "

$ws.Range("D30").Value = "'1.468"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").Value = "'0.05724"
$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").Value = "'4.124"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").Value = "'4.047"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").Value = "'0.7082"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D38").Value = "'2.776"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'0.01792"
$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("D40").Value = "'1.216.98"
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("D41").Value = "'6.500"
$ws.Range("E41").Value = "  +5.06%  "

$ws.Range("D42").Value = "'0.9074"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'2.016.81"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").Value = "'101.85"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "'66.35"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").Value = "'7.126"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "'8.977"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("D51").Value = "'1.681"
$ws.Range("E51").Value = "  -0.33%  "
